$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pl_mw results for the 380 kV case (Case_5_26): rows 2-25, columns B-F and J-M

# Row 2
$ws.Range("B2").Value = 1.770638372978169
$ws.Range("C2").Value = 0.01311926495755245
$ws.Range("D2").Value = 0.03223974428605914
$ws.Range("E2").Value = 0.06383820745122026
$ws.Range("F2").Value = 7.662150779775345
$ws.Range("J2").Value = 0.2588048803499703
$ws.Range("K2").Value = 1.23689735626985
$ws.Range("L2").Value = 0.2246228339345819
$ws.Range("M2").Value = 0.3682373811667823

# Row 3
$ws.Range("B3").Value = 1.76829961054716
$ws.Range("C3").Value = 0.01138622289428781
$ws.Range("D3").Value = 0.02831459372056599
$ws.Range("E3").Value = 0.06421076739119824
$ws.Range("F3").Value = 7.468097427214133
$ws.Range("J3").Value = 0.2557179444047435
$ws.Range("K3").Value = 1.231879661895704
$ws.Range("L3").Value = 0.2269173297750591
$ws.Range("M3").Value = 0.3693921511879914

# Row 4
$ws.Range("B4").Value = 1.768133792704134
$ws.Range("C4").Value = 0.01033676375450909
$ws.Range("D4").Value = 0.02589476758404174
$ws.Range("E4").Value = 0.06445473036819926
$ws.Range("F4").Value = 7.349342984866155
$ws.Range("J4").Value = 0.2538144995235854
$ws.Range("K4").Value = 1.229838432415718
$ws.Range("L4").Value = 0.2284582534628612
$ws.Range("M4").Value = 0.3703540916250851

# Row 5
$ws.Range("B5").Value = 1.768385646666843
$ws.Range("C5").Value = 0.009912679052497708
$ws.Range("D5").Value = 0.02490607137529111
$ws.Range("E5").Value = 0.06455798277132851
$ws.Range("F5").Value = 7.301045931885881
$ws.Range("J5").Value = 0.2530367374710352
$ws.Range("K5").Value = 1.229267871344064
$ws.Range("L5").Value = 0.2291194510675361
$ws.Range("M5").Value = 0.3708097158614514

# Row 6
$ws.Range("B6").Value = 1.768446760421568
$ws.Range("C6").Value = 0.009842473273231178
$ws.Range("D6").Value = 0.02474173812076685
$ws.Range("E6").Value = 0.0645753597308194
$ws.Range("F6").Value = 7.29303197494383
$ws.Range("J6").Value = 0.2529074623055365
$ws.Range("K6").Value = 1.229188904546945
$ws.Range("L6").Value = 0.2292312524355289
$ws.Range("M6").Value = 0.3708892156879564

# Row 7
$ws.Range("B7").Value = 1.768135895918249
$ws.Range("C7").Value = 0.01033103004449742
$ws.Range("D7").Value = 0.02588144438318096
$ws.Range("E7").Value = 0.06445610732175222
$ws.Range("F7").Value = 7.348691247812411
$ws.Range("J7").Value = 0.2538040189133071
$ws.Range("K7").Value = 1.229829680041846
$ws.Range("L7").Value = 0.2284670358861973
$ws.Range("M7").Value = 0.3703599786753777

# Row 8
$ws.Range("B8").Value = 1.769568343775347
$ws.Range("C8").Value = 0.01251861243876817
$ws.Range("D8").Value = 0.03088829371017709
$ws.Range("E8").Value = 0.06396351590999583
$ws.Range("F8").Value = 7.595157307588323
$ws.Range("J8").Value = 0.257742134804225
$ws.Range("K8").Value = 1.234951360355979
$ws.Range("L8").Value = 0.2253865852761905
$ws.Range("M8").Value = 0.3685830714677607

# Row 9
$ws.Range("B9").Value = 1.782458746692413
$ws.Range("C9").Value = 0.01692939683604067
$ws.Range("D9").Value = 0.04063595127189501
$ws.Range("E9").Value = 0.06311773911260854
$ws.Range("F9").Value = 8.081771553619518
$ws.Range("J9").Value = 0.2654042491395501
$ws.Range("K9").Value = 1.25325567135215
$ws.Range("L9").Value = 0.2203921392919597
$ws.Range("M9").Value = 0.3671045329615481

# Row 10
$ws.Range("B10").Value = 1.798085353570201
$ws.Range("C10").Value = 0.02025071311200577
$ws.Range("D10").Value = 0.04776429543628069
$ws.Range("E10").Value = 0.06256895009657804
$ws.Range("F10").Value = 8.4415528705741
$ws.Range("J10").Value = 0.2710018749575198
$ws.Range("K10").Value = 1.271761115055767
$ws.Range("L10").Value = 0.2173582705303687
$ws.Range("M10").Value = 0.3672410474560159

# Row 11
$ws.Range("B11").Value = 1.806533849416041
$ws.Range("C11").Value = 0.02178072871076608
$ws.Range("D11").Value = 0.05100210241022296
$ws.Range("E11").Value = 0.06233491532617963
$ws.Range("F11").Value = 8.605776620533959
$ws.Range("J11").Value = 0.2735426346901448
$ws.Range("K11").Value = 1.281283015418865
$ws.Range("L11").Value = 0.216115626335025
$ws.Range("M11").Value = 0.3675687559826564

# Row 12
$ws.Range("B12").Value = 1.809925905343079
$ws.Range("C12").Value = 0.02236297684790145
$ws.Range("D12").Value = 0.05222763200308123
$ws.Range("E12").Value = 0.0622485264713255
$ws.Range("F12").Value = 8.668048006292452
$ws.Range("J12").Value = 0.2745040272891117
$ws.Range("K12").Value = 1.285047772409541
$ws.Range("L12").Value = 0.2156648027483925
$ws.Range("M12").Value = 0.3677310434028129

# Row 13
$ws.Range("B13").Value = 1.809186788744427
$ws.Range("C13").Value = 0.02223745008197398
$ws.Range("D13").Value = 0.05196371470390204
$ws.Range("E13").Value = 0.06226703261111677
$ws.Range("F13").Value = 8.654632986926231
$ws.Range("J13").Value = 0.2742970058709275
$ws.Range("K13").Value = 1.284229887052902
$ws.Range("L13").Value = 0.2157610182575525
$ws.Range("M13").Value = 0.3676943933480423

# Row 14
$ws.Range("B14").Value = 1.806809051271244
$ws.Range("C14").Value = 0.02182857255445469
$ws.Range("D14").Value = 0.05110293778835739
$ws.Range("E14").Value = 0.06232776332654666
$ws.Range("F14").Value = 8.610898034335946
$ws.Range("J14").Value = 0.2736217434170598
$ws.Range("K14").Value = 1.281589555407578
$ws.Range("L14").Value = 0.2160781413823329
$ws.Range("M14").Value = 0.3675813420615803

# Row 15
$ws.Range("B15").Value = 1.805377729740343
$ws.Range("C15").Value = 0.02157849967944969
$ws.Range("D15").Value = 0.05057561888042983
$ws.Range("E15").Value = 0.06236525341615629
$ws.Range("F15").Value = 8.584120073638132
$ws.Range("J15").Value = 0.2732080316884193
$ws.Range("K15").Value = 1.279992994633318
$ws.Range("L15").Value = 0.2162749581443322
$ws.Range("M15").Value = 0.3675170684862259

# Row 16
$ws.Range("B16").Value = 1.797560196240141
$ws.Range("C16").Value = 0.02015111763634536
$ws.Range("D16").Value = 0.04755261098614483
$ws.Range("E16").Value = 0.06258455812071029
$ws.Range("F16").Value = 8.430831956774796
$ws.Range("J16").Value = 0.2708357231668472
$ws.Range("K16").Value = 1.271161071702835
$ws.Range("L16").Value = 0.2174422440389208
$ws.Range("M16").Value = 0.367224975339834

# Row 17
$ws.Range("B17").Value = 1.793107648253823
$ws.Range("C17").Value = 0.01928044539636176
$ws.Range("D17").Value = 0.04569694100281652
$ws.Range("E17").Value = 0.06272308596978249
$ws.Range("F17").Value = 8.336939738166336
$ws.Range("J17").Value = 0.2693789978016383
$ws.Range("K17").Value = 1.266025869729731
$ws.Range("L17").Value = 0.2181935255838994
$ws.Range("M17").Value = 0.3671138112972123

# Row 18
$ws.Range("B18").Value = 1.790672767517009
$ws.Range("C18").Value = 0.01878145171895085
$ws.Range("D18").Value = 0.04462913738497321
$ws.Range("E18").Value = 0.06280423360120091
$ws.Range("F18").Value = 8.282987687559597
$ws.Range("J18").Value = 0.268540593836029
$ws.Range("K18").Value = 1.263176101569115
$ws.Range("L18").Value = 0.2186385846832337
$ws.Range("M18").Value = 0.3670748729065636

# Row 19
$ws.Range("B19").Value = 1.789870014279728
$ws.Range("C19").Value = 0.01861280571722546
$ws.Range("D19").Value = 0.04426751195423151
$ws.Range("E19").Value = 0.06283196158678805
$ws.Range("F19").Value = 8.264729329039227
$ws.Range("J19").Value = 0.2682566307675671
$ws.Range("K19").Value = 1.262229048526848
$ws.Range("L19").Value = 0.2187914977585521
$ws.Range("M19").Value = 0.3670659829731804

# Row 20
$ws.Range("B20").Value = 1.793568577179968
$ws.Range("C20").Value = 0.01937294345663076
$ws.Range("D20").Value = 0.04589452794135696
$ws.Range("E20").Value = 0.0627081873690285
$ws.Range("F20").Value = 8.346929291914961
$ws.Range("J20").Value = 0.2695341234583708
$ws.Range("K20").Value = 1.266561769006785
$ws.Range("L20").Value = 0.2181122111777256
$ws.Range("M20").Value = 0.3671230575053173

# Row 21
$ws.Range("B21").Value = 1.807502216987132
$ws.Range("C21").Value = 0.02194859114695191
$ws.Range("D21").Value = 0.0513557827799076
$ws.Range("E21").Value = 0.06230986466417443
$ws.Range("F21").Value = 8.623741760012706
$ws.Range("J21").Value = 0.2738201038255283
$ws.Range("K21").Value = 1.282360766500346
$ws.Range("L21").Value = 0.2159844591649502
$ws.Range("M21").Value = 0.3676135114715855

# Row 22
$ws.Range("B22").Value = 1.817732473885286
$ws.Range("C22").Value = 0.02364867535325743
$ws.Range("D22").Value = 0.05492184768695552
$ws.Range("E22").Value = 0.0620625606677101
$ws.Range("F22").Value = 8.805142461588105
$ws.Range("J22").Value = 0.2766169718090339
$ws.Range("K22").Value = 1.293613369808838
$ws.Range("L22").Value = 0.2147088897313907
$ws.Range("M22").Value = 0.368156660373792

# Row 23
$ws.Range("B23").Value = 1.812169519644755
$ws.Range("C23").Value = 0.02273973809677443
$ws.Range("D23").Value = 0.05301881345287995
$ws.Range("E23").Value = 0.06219336313935209
$ws.Range("F23").Value = 8.708279776794427
$ws.Range("J23").Value = 0.275124597525938
$ws.Range("K23").Value = 1.287522712131221
$ws.Range("L23").Value = 0.2153791687734241
$ws.Range("M23").Value = 0.3678464030391417

# Row 24
$ws.Range("B24").Value = 1.793359802020206
$ws.Range("C24").Value = 0.01933112021080774
$ws.Range("D24").Value = 0.04580520187117543
$ws.Range("E24").Value = 0.0627149183343616
$ws.Range("F24").Value = 8.342412928184558
$ws.Range("J24").Value = 0.2694639940008301
$ws.Range("K24").Value = 1.266319169603491
$ws.Range("L24").Value = 0.2181489324968737
$ws.Range("M24").Value = 0.3671187995112781

# Row 25
$ws.Range("B25").Value = 1.777891274929857
$ws.Range("C25").Value = 0.01572249727995967
$ws.Range("D25").Value = 0.03800550785004475
$ws.Range("E25").Value = 0.0633337461647776
$ws.Range("F25").Value = 7.94975031032422
$ws.Range("J25").Value = 0.2633373910600127
$ws.Range("K25").Value = 1.247417350917459
$ws.Range("L25").Value = 0.2216314970005442
$ws.Range("M25").Value = 0.3672897833950692

Write-Output "applied 380 kV case updates"
